$wb = $excel.ActiveWorkbook

# Update the "zh-cn" sheet handoff/handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-08 11:13:38"
$wsZhCn.Range("G2").Value = "2016-01-08 11:14:20"

# Update the "de-de" sheet handoff/handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-08 11:13:47"
$wsDeDe.Range("G2").Value = "2016-01-08 11:14:36"
